$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the value to be stored as text (matches source cells, which are
    # inline/shared strings even when the text looks numeric), then clear the
    # temporary Text number-format so no stray style is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "70.123.03"
Set-TextValue $ws.Range("E2") "  +0.64%  "

Set-TextValue $ws.Range("D3") "3.585.56"
Set-TextValue $ws.Range("E3") "  +0.15%  "

Set-TextValue $ws.Range("E4") "  +0.21%  "

Set-TextValue $ws.Range("D5") "574.80"
Set-TextValue $ws.Range("E5") "  -3.17%  "

Set-TextValue $ws.Range("D6") "190.29"
Set-TextValue $ws.Range("E6") "  -1.17%  "

Set-TextValue $ws.Range("D7") "0.632"
Set-TextValue $ws.Range("E7") "  -1.60%  "

Set-TextValue $ws.Range("D8") "3.584.71"
Set-TextValue $ws.Range("E8") "  +0.32%  "

Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  +0.00%  "

Set-TextValue $ws.Range("E10") "  -3.79%  "

Set-TextValue $ws.Range("D11") "0.659"
Set-TextValue $ws.Range("E11") "  -0.35%  "

Set-TextValue $ws.Range("D12") "56.43"
Set-TextValue $ws.Range("E12") "  -3.14%  "

Set-TextValue $ws.Range("D13") "0.0000299"
Set-TextValue $ws.Range("E13") "  +2.31%  "

Set-TextValue $ws.Range("D14") "9.77"
Set-TextValue $ws.Range("E14") "  +0.57%  "

Set-TextValue $ws.Range("D15") "4.153.99"
Set-TextValue $ws.Range("E15") "  +0.35%  "

Set-TextValue $ws.Range("D16") "20.09"
Set-TextValue $ws.Range("E16") "  +3.92%  "

Set-TextValue $ws.Range("D17") "3.579.63"
Set-TextValue $ws.Range("E17") "  +0.30%  "

Set-TextValue $ws.Range("D18") "69.949.98"
Set-TextValue $ws.Range("E18") "  +0.83%  "

Set-TextValue $ws.Range("D19") "12.51"
Set-TextValue $ws.Range("E19") "  +0.27%  "

Set-TextValue $ws.Range("E20") "  +1.00%  "

Set-TextValue $ws.Range("E21") "  -0.87%  "

Set-TextValue $ws.Range("B22") "BitcoinCash"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D22") "472.57"
Set-TextValue $ws.Range("E22") "  -5.82%  "

Set-TextValue $ws.Range("B23") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D23") "19.60"
Set-TextValue $ws.Range("E23") "  +13.86%  "

Set-TextValue $ws.Range("E24") "  -6.68%  "

Set-TextValue $ws.Range("D25") "4.33"
Set-TextValue $ws.Range("E25") "  -2.63%  "

Set-TextValue $ws.Range("D26") "88.49"
Set-TextValue $ws.Range("E26") "  -3.05%  "

Set-TextValue $ws.Range("E27") "  -0.75%  "

Set-TextValue $ws.Range("D28") "11.10"
Set-TextValue $ws.Range("E28") "  -0.86%  "

Set-TextValue $ws.Range("D29") "9.25"
Set-TextValue $ws.Range("E29") "  -1.15%  "

Set-TextValue $ws.Range("D30") "7.78"
Set-TextValue $ws.Range("E30") "  +3.09%  "

Set-TextValue $ws.Range("D31") "32.08"
Set-TextValue $ws.Range("E31") "  -0.22%  "

Set-TextValue $ws.Range("D32") "0.120"
Set-TextValue $ws.Range("E32") "  +3.88%  "

Set-TextValue $ws.Range("D33") "12.09"
Set-TextValue $ws.Range("E33") "  -0.19%  "

Set-TextValue $ws.Range("D34") "66.25"
Set-TextValue $ws.Range("E34") "  +1.18%  "

Set-TextValue $ws.Range("D35") "589.13"
Set-TextValue $ws.Range("E35") "  -4.39%  "

Set-TextValue $ws.Range("D36") "39.52"
Set-TextValue $ws.Range("E36") "  +3.64%  "

Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  -0.01%  "

Set-TextValue $ws.Range("D38") "0.0₃0804"
Set-TextValue $ws.Range("E38") "  -4.19%  "

Set-TextValue $ws.Range("D39") "0.399"
Set-TextValue $ws.Range("E39") "  +0.03%  "

Set-TextValue $ws.Range("D40") "0.143"
Set-TextValue $ws.Range("E40") "  -3.41%  "

Set-TextValue $ws.Range("D41") "3.53"
Set-TextValue $ws.Range("E41") "  -2.58%  "

Set-TextValue $ws.Range("D42") "2.92"
Set-TextValue $ws.Range("E42") "  +7.68%  "

Set-TextValue $ws.Range("D43") "3.228.49"
Set-TextValue $ws.Range("E43") "  -3.05%  "

Set-TextValue $ws.Range("D44") "3.13"
Set-TextValue $ws.Range("E44") "  +7.99%  "

Set-TextValue $ws.Range("D45") "3.10"
Set-TextValue $ws.Range("E45") "  -0.58%  "

Set-TextValue $ws.Range("D46") "0.0444"
Set-TextValue $ws.Range("E46") "  +0.47%  "

Set-TextValue $ws.Range("D47") "9.61"
Set-TextValue $ws.Range("E47") "  +5.08%  "

Set-TextValue $ws.Range("D48") "3.35"
Set-TextValue $ws.Range("E48") "  +1.15%  "

Set-TextValue $ws.Range("E49") "  -0.65%  "

Set-TextValue $ws.Range("D50") "1.00"
Set-TextValue $ws.Range("E50") "  +0.41%  "

Set-TextValue $ws.Range("D51") "3.16"
Set-TextValue $ws.Range("E51") "  -2.51%  "
